# Replace embedded line-breaks inside a handful of vaccine name / brand cells
# with plain spaces, on the two "Influenza" sheets.
#
# "Pediatric Influenza Vaccine " sheet
$ws3 = $excel.ActiveWorkbook.Worksheets.Item("Pediatric Influenza Vaccine ")
$ws3.Range("B3").Value  = "Fluzone Pediatric dose No Preservative"
$ws3.Range("B6").Value  = "Fluarix Preservative-Free"
$ws3.Range("B9").Value  = "FluMist No Preservative"
$ws3.Range("B10").Value = "Afluria No Preservative"
$ws3.Range("H10").Value = "Merck (CSL product)"

# "Adult Influenza Vaccine " sheet
$ws4 = $excel.ActiveWorkbook.Worksheets.Item("Adult Influenza Vaccine ")
$ws4.Range("B5").Value  = "Agriflu No Preservative"
$ws4.Range("B7").Value  = "Fluvirin Preservative-free"
$ws4.Range("B10").Value = "Flumist No Preservative"
